$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers I1 "I0" and J1 "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the I0 / IF values for each data row (2-44)
$data = @(
    @(2, 10, 10),
    @(3, 7, 8),
    @(4, 5, 6),
    @(5, 7, 8),
    @(6, 8, 8),
    @(7, 8, 8),
    @(8, 5, 5),
    @(9, 7, 7),
    @(10, 6, 7),
    @(11, 8, 9),
    @(12, 6, 7),
    @(13, 4, 6),
    @(14, 7, 7),
    @(15, 6, 6),
    @(16, 8, 8),
    @(17, 7, 8),
    @(18, 7, 7),
    @(19, 5, 7),
    @(20, 8, 8),
    @(21, 7, 7),
    @(22, 7, 7),
    @(23, 8, 8),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 9, 9),
    @(27, 9, 9),
    @(28, 9, 9),
    @(29, 9, 9),
    @(30, 9, 9),
    @(31, 9, 9),
    @(32, 8, 9),
    @(33, 8, 9),
    @(34, 9, 9),
    @(35, 9, 9),
    @(36, 7, 7),
    @(37, 9, 9),
    @(38, 8, 10),
    @(39, 9, 9),
    @(40, 9, 9),
    @(41, 6, 7),
    @(42, 1, 3),
    @(43, 5, 5),
    @(44, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
